# Applies the "x2 is the best" commit: pastes 100 freshly-recorded solver
# samples (time in column A, move count in column B) into the analytics
# sheet "NyanyanFunc_-1.9-2事前計算_3_max_sd" (the workbook's first / active
# tab, physical part xl/worksheets/sheet1.xml). All of the AVERAGE / STDEV /
# MAX / MEDIAN / MIN / COUNTIFS formulas on that sheet recalc automatically
# once the data lands. The only other change is a remembered-selection
# (view-state) tweak on the second tab.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (tab 1, tabSelected, ActiveSheet): "NyanyanFunc_-1.9-2事前計算_3_max_sd"
$ws = $wb.Worksheets.Item(1)

$arr = New-Object 'object[,]' 100,2
$arr[0,0] = 5.0672748089999997
$arr[0,1] = 58
$arr[1,0] = 0.287303209
$arr[1,1] = 51
$arr[2,0] = 6.5067868229999997
$arr[2,1] = 61
$arr[3,0] = 1.440039635
$arr[3,1] = 58
$arr[4,0] = 0.66227531399999995
$arr[4,1] = 53
$arr[5,0] = 0.92621779400000004
$arr[5,1] = 53
$arr[6,0] = 2.8318591120000001
$arr[6,1] = 59
$arr[7,0] = 1.5928866859999999
$arr[7,1] = 56
$arr[8,0] = 5.0022780899999999
$arr[8,1] = 58
$arr[9,0] = 5.3004183769999997
$arr[9,1] = 55
$arr[10,0] = 6.6854703430000004
$arr[10,1] = 57
$arr[11,0] = 2.9165472979999998
$arr[11,1] = 56
$arr[12,0] = 1.670004606
$arr[12,1] = 55
$arr[13,0] = 8.2502539160000001
$arr[13,1] = 58
$arr[14,0] = 11.188903809999999
$arr[14,1] = 61
$arr[15,0] = 0.74304985999999995
$arr[15,1] = 53
$arr[16,0] = 0.71809864000000001
$arr[16,1] = 54
$arr[17,0] = 8.778671503
$arr[17,1] = 61
$arr[18,0] = 4.8911504749999999
$arr[18,1] = 58
$arr[19,0] = 2.9650721550000001
$arr[19,1] = 58
$arr[20,0] = 1.9836964610000001
$arr[20,1] = 53
$arr[21,0] = 3.5445199010000001
$arr[21,1] = 54
$arr[22,0] = 2.8025181290000001
$arr[22,1] = 57
$arr[23,0] = 4.4311487669999998
$arr[23,1] = 58
$arr[24,0] = 4.1358075139999997
$arr[24,1] = 55
$arr[25,0] = 7.7010188099999999
$arr[25,1] = 58
$arr[26,0] = 3.8961384300000002
$arr[26,1] = 61
$arr[27,0] = 4.2017691140000002
$arr[27,1] = 58
$arr[28,0] = 1.5299046039999999
$arr[28,1] = 55
$arr[29,0] = 11.08447003
$arr[29,1] = 60
$arr[30,0] = 1.7592420580000001
$arr[30,1] = 54
$arr[31,0] = 3.7530398370000002
$arr[31,1] = 58
$arr[32,0] = 3.5784304140000001
$arr[32,1] = 56
$arr[33,0] = 5.8275206089999996
$arr[33,1] = 56
$arr[34,0] = 0.94942855800000003
$arr[34,1] = 52
$arr[35,0] = 2.0784864430000001
$arr[35,1] = 56
$arr[36,0] = 2.5875957010000001
$arr[36,1] = 59
$arr[37,0] = 1.6236579419999999
$arr[37,1] = 58
$arr[38,0] = 1.0193178650000001
$arr[38,1] = 52
$arr[39,0] = 1.732247353
$arr[39,1] = 54
$arr[40,0] = 3.834183216
$arr[40,1] = 59
$arr[41,0] = 2.8641924859999999
$arr[41,1] = 57
$arr[42,0] = 13.00402379
$arr[42,1] = 57
$arr[43,0] = 3.7727608680000002
$arr[43,1] = 58
$arr[44,0] = 8.5419690609999996
$arr[44,1] = 58
$arr[45,0] = 0.66023397399999995
$arr[45,1] = 54
$arr[46,0] = 4.0819411280000004
$arr[46,1] = 59
$arr[47,0] = 9.46049118
$arr[47,1] = 61
$arr[48,0] = 7.0185046199999999
$arr[48,1] = 57
$arr[49,0] = 1.6422502990000001
$arr[49,1] = 55
$arr[50,0] = 3.1352026460000002
$arr[50,1] = 58
$arr[51,0] = 1.2506206040000001
$arr[51,1] = 57
$arr[52,0] = 3.9224591260000001
$arr[52,1] = 59
$arr[53,0] = 7.8286769390000002
$arr[53,1] = 58
$arr[54,0] = 4.3092212679999999
$arr[54,1] = 56
$arr[55,0] = 7.0957560539999998
$arr[55,1] = 58
$arr[56,0] = 0.31914639500000003
$arr[56,1] = 52
$arr[57,0] = 3.4268589020000002
$arr[57,1] = 55
$arr[58,0] = 3.2952854629999999
$arr[58,1] = 58
$arr[59,0] = 2.42135334
$arr[59,1] = 59
$arr[60,0] = 1.662703037
$arr[60,1] = 58
$arr[61,0] = 3.982432604
$arr[61,1] = 59
$arr[62,0] = 4.0896730420000003
$arr[62,1] = 54
$arr[63,0] = 0.23536992100000001
$arr[63,1] = 52
$arr[64,0] = 5.5029833320000003
$arr[64,1] = 59
$arr[65,0] = 5.9727659229999999
$arr[65,1] = 59
$arr[66,0] = 4.0769577029999997
$arr[66,1] = 59
$arr[67,0] = 2.5970783229999999
$arr[67,1] = 55
$arr[68,0] = 8.9266226290000006
$arr[68,1] = 60
$arr[69,0] = 11.640734670000001
$arr[69,1] = 55
$arr[70,0] = 1.3334331509999999
$arr[70,1] = 56
$arr[71,0] = 3.6005725860000002
$arr[71,1] = 60
$arr[72,0] = 4.4960887429999996
$arr[72,1] = 60
$arr[73,0] = 0.41588640199999999
$arr[73,1] = 50
$arr[74,0] = 1.0907924179999999
$arr[74,1] = 58
$arr[75,0] = 13.147104499999999
$arr[75,1] = 55
$arr[76,0] = 4.2637541289999996
$arr[76,1] = 60
$arr[77,0] = 5.8724448679999997
$arr[77,1] = 61
$arr[78,0] = 4.4637324810000001
$arr[78,1] = 58
$arr[79,0] = 1.2601177690000001
$arr[79,1] = 56
$arr[80,0] = 6.7468934059999999
$arr[80,1] = 58
$arr[81,0] = 4.8749299050000001
$arr[81,1] = 56
$arr[82,0] = 1.568384886
$arr[82,1] = 53
$arr[83,0] = 1.463377476
$arr[83,1] = 56
$arr[84,0] = 14.11067486
$arr[84,1] = 56
$arr[85,0] = 1.2223513130000001
$arr[85,1] = 54
$arr[86,0] = 2.9281678200000001
$arr[86,1] = 59
$arr[87,0] = 1.8511536120000001
$arr[87,1] = 55
$arr[88,0] = 3.5030717849999999
$arr[88,1] = 60
$arr[89,0] = 1.0013461109999999
$arr[89,1] = 55
$arr[90,0] = 5.432256937
$arr[90,1] = 56
$arr[91,0] = 1.702488899
$arr[91,1] = 57
$arr[92,0] = 6.574088573
$arr[92,1] = 60
$arr[93,0] = 1.4726161959999999
$arr[93,1] = 53
$arr[94,0] = 8.6567151550000005
$arr[94,1] = 56
$arr[95,0] = 9.2558860779999996
$arr[95,1] = 58
$arr[96,0] = 6.6203200820000001
$arr[96,1] = 60
$arr[97,0] = 3.8141021730000002
$arr[97,1] = 59
$arr[98,0] = 4.4641282560000004
$arr[98,1] = 56
$arr[99,0] = 3.725548506
$arr[99,1] = 57
$ws.Range("A2:B101").Value = $arr

# --- Sheet 2 (tab 2, not active): "NyanyanFunc_-1.9-2事前計算_2_max_sd"
# Only its cached selection moves (P28 -> D30:D31); no data changes there.
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("D30:D31").Select()

# Re-select sheet 1 last so it ends up the active tab again (tabSelected="1")
# with its own selection left on G16, matching the saved view state.
$ws.Range("G16").Select()
